$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.115.54'
$ws.Range("E2").Value = '  +0.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.666.37'

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.82'
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5205'
$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("E7").Value = '  -0.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2601'
$ws.Range("E8").Value = '  -2.61%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06320'
$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.05'
$ws.Range("E10").Value = '  -0.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07529'
$ws.Range("E11").Value = '  -0.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.680.39'
$ws.Range("E12").Value = '  -0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.410'
$ws.Range("E13").Value = '  -1.89%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5420'
$ws.Range("E14").Value = '  -4.16%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000007997'
$ws.Range("E15").Value = '  -1.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.31'
$ws.Range("E16").Value = '  +1.31%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.168.68'
$ws.Range("E17").Value = '  +0.22%  '

$ws.Range("E18").Value = '  -0.23%  '

$ws.Range("E19").Value = '  -2.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '187.02'
$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.25'
$ws.Range("E21").Value = '  -3.17%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.229'
$ws.Range("E22").Value = '  +0.78%  '

$ws.Range("E23").Value = '  -0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.88'
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1233'
$ws.Range("E25").Value = '  -1.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.453'
$ws.Range("E26").Value = '  -2.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.73'
$ws.Range("E27").Value = '  -0.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06267'
$ws.Range("E28").Value = '  -1.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.365'
$ws.Range("E29").Value = '  +0.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.278'
$ws.Range("E30").Value = '  -0.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.494'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.405'
$ws.Range("E32").Value = '  -3.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.640'
$ws.Range("E33").Value = '  -1.06%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9990'
$ws.Range("E34").Value = '  -0.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.110.05'
$ws.Range("E38").Value = '  +1.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01608'
$ws.Range("E39").Value = '  -0.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.055'
$ws.Range("E40").Value = '  -1.38%  '

$ws.Range("E41").Value = '  -0.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.15%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.58'
$ws.Range("E43").Value = '  +0.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.817.92'
$ws.Range("E44").Value = '  -0.47%  '

$ws.Range("E45").Value = '  -0.23%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.29'
$ws.Range("E46").Value = '  -2.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.004'
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.053'
$ws.Range("E48").Value = '  +1.25%  '

$ws.Range("E49").Value = '  -0.12%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4234'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.879'
$ws.Range("E51").Value = '  -0.75%  '

# Row 35 <-> Row 36 swap (HuobiToken <-> ImmutableX)
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.5992'
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.393'
$ws.Range("E36").Value = '  -0.92%  '
